$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.321358333333333
$ws.Range("N2").Value = 3.964075
$ws.Range("O2").Value = 0.06904315418552966
$ws.Range("P2").Value = 0.06904315418552966
$ws.Range("Q2").Value = 0.2135315066666666
$ws.Range("R2").Value = 1.92178356
$ws.Range("S2").Value = 0.06904315418552966
$ws.Range("T2").Value = 0.06904315418552966

# Row 3
$ws.Range("O3").Value = 0.4558096119837698
$ws.Range("P3").Value = 0.4558096119837698
$ws.Range("S3").Value = 0.4558096119837698
$ws.Range("T3").Value = 0.4558096119837698

# Row 4
$ws.Range("M4").Value = 9.093439666666667
$ws.Range("N4").Value = 27.280319
$ws.Range("O4").Value = 0.4751472338307006
$ws.Range("P4").Value = 0.4751472338307005
$ws.Range("Q4").Value = 1.469499850133333
$ws.Range("R4").Value = 13.2254986512
$ws.Range("S4").Value = 0.4751472338307006
$ws.Range("T4").Value = 0.4751472338307005
